# Update Frambuesa (raspberry) price records in the Lo Valledor wholesale
# market sheet: rewrite rows 259-261 with new weekly figures and append two
# new rows (263, 264) carrying the data that used to live in row 262, plus a
# brand-new row (262) with another Linares "Primera" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Set-Row {
    param([int]$Row, [double]$Fecha, [string]$Calidad, [double]$Volumen, [double]$PrecioMinimo, [double]$PrecioMaximo, [double]$PrecioPromedio, [string]$Origen, [double]$PrecioKg)

    $ws.Range("A$Row").Value = 6
    $ws.Range("B$Row").Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Range("C$Row").Value = "Metropolitana"

    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("D$Row").NumberFormat = $dateFormat

    $ws.Range("E$Row").Value = 13
    $ws.Range("F$Row").Value = "Fruta"
    $ws.Range("G$Row").Value = 100101
    $ws.Range("H$Row").Value = "Berries"
    $ws.Range("I$Row").Value = 100101004
    $ws.Range("J$Row").Value = "Frambuesa"
    $ws.Range("K$Row").Value = "Sin especificar"
    $ws.Range("L$Row").Value = $Calidad
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $PrecioMinimo
    $ws.Range("O$Row").Value = $PrecioMaximo
    $ws.Range("P$Row").Value = $PrecioPromedio
    $ws.Range("Q$Row").Value = "$/bandeja 2 kilos"
    $ws.Range("R$Row").Value = $Origen
    $ws.Range("S$Row").Value = $PrecioKg
    $ws.Range("T$Row").Value = 2
}

# Row 259: 2022-05-09 -> 2023-01-13, Primera -> Especial
Set-Row 259 44939 "Especial" 200 8000 8000 8000 "Provincia de Curicó" 4000

# Row 260: 2021-01-13 -> 2023-01-13, Primera -> Especial, Linares -> Región del Maule
Set-Row 260 44939 "Especial" 250 8000 8000 8000 "Región del Maule" 4000

# Row 261: 2021-01-13 -> 2022-05-09, Segunda -> Primera, Linares -> Curicó
Set-Row 261 44690 "Primera" 75 12000 12000 12000 "Provincia de Curicó" 6000

# Row 262 (new): 2021-01-13, Primera, Linares
Set-Row 262 44209 "Primera" 185 6000 6000 6000 "Provincia de Linares" 3000

# Row 263 (new): 2021-01-13, Segunda, Linares
Set-Row 263 44209 "Segunda" 120 4000 4000 4000 "Provincia de Linares" 2000

# Row 264 (new): carries forward the record that used to be row 262 (2021-05-19, Primera, Curicó)
Set-Row 264 44335 "Primera" 75 10000 10000 10000 "Provincia de Curicó" 5000
